$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "zMaxShap"
$ws.Range("B9").Value = "asfgfkjlahfs"
$ws.Range("C9").Value = "Sman95@gmail.com"

$ws.Range("A10").Value = "sfdafg"
$ws.Range("B10").Value = "asdasfg"
$ws.Range("C10").Value = "Sman95@gmail.com"
